# Auto-generated script applying the Famfrit_Profits data refresh diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H33" = 597.7273
    "I33" = 505.88235
    "K33" = 505.88235
    "M33" = -276.88235
    "H88" = 2700.077
    "I88" = 2251.5
    "J88" = 2781.6365
    "K88" = 2251.5
    "L88" = 2781.6365
    "M88" = -1845.5
    "N88" = -3593.6365
    "H91" = 2700.077
    "I91" = 2251.5
    "J91" = 2781.6365
    "K91" = 2251.5
    "L91" = 2781.6365
    "M91" = -847.5
    "N91" = -5589.636500000001
    "H98" = 634.381
    "I98" = 634.381
    "K98" = 634.381
    "M98" = 863.619
    "H101" = 618.2143
    "I101" = 674.125
    "K101" = 2022.375
    "M101" = -400.375
    "H122" = 634.381
    "I122" = 634.381
    "K122" = 1903.143
    "M122" = 546.857
    "H132" = 3240.561
    "I132" = 3564.853
    "J132" = 1665.4286
    "K132" = 10694.559
    "L132" = 4996.2858
    "M132" = -8164.559000000001
    "N132" = -10056.2858
    "H135" = 1134.3
    "I135" = 924.3333
    "J135" = 1449.25
    "K135" = 8318.9997
    "L135" = 13043.25
    "M135" = -5783.9997
    "N135" = -18113.25
    "H137" = 15899.348
    "I137" = 22141.188
    "K137" = 66423.564
    "M137" = -63873.564
    "H141" = 1722.4615
    "I141" = 1308.3
    "K141" = 3924.9
    "M141" = 1255.1
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H29" = 500
    "I29" = 500
    "J29" = 0
    "K29" = 500
    "L29" = 0
    "M29" = -192
    "H61" = 2982.5862
    "I61" = 2688.087
    "K61" = 2688.087
    "M61" = -2476.087
    "H74" = 48862.72
    "I74" = 52698.61
    "J74" = 4750
    "K74" = 52698.61
    "L74" = 4750
    "M74" = -51824.61
    "N74" = -6498
    "H77" = 48862.72
    "I77" = 52698.61
    "J77" = 4750
    "K77" = 263493.05
    "L77" = 23750
    "M77" = -259125.05
    "N77" = -32486
    "H97" = 1050.2
    "I97" = 933.619
    "K97" = 933.619
    "M97" = -437.619
    "H122" = 1484.317
    "I122" = 1229.9459
    "K122" = 3689.8377
    "M122" = -1239.8377
    "H133" = 104065
    "J133" = 104065
    "L133" = 104065
    "N133" = -109125
    "H136" = 2982.5862
    "I136" = 2688.087
    "K136" = 8064.261
    "M136" = -5514.261
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H86" = 25090.77
    "I86" = 11440.5
    "K86" = 11440.5
    "M86" = -10317.5
    "H89" = 25090.77
    "I89" = 11440.5
    "K89" = 57202.5
    "M89" = -51586.5
    "H105" = 13994.08
    "I105" = 26141.889
    "J105" = 7160.9375
    "K105" = 26141.889
    "L105" = 7160.9375
    "M105" = -24394.889
    "N105" = -10654.9375
    "H106" = 234980.67
    "J106" = 234980.67
    "L106" = 234980.67
    "N106" = -237504.67
    "H132" = 139780
    "J132" = 139780
    "L132" = 139780
    "N132" = -149900
    "H134" = 1988.6666
    "I134" = 987.375
    "K134" = 2962.125
    "M134" = -427.125
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 4152.3335
    "I31" = 1895.1666
    "J31" = 6409.5
    "K31" = 1895.1666
    "L31" = 6409.5
    "M31" = -1600.1666
    "N31" = -6999.5
    "H34" = 4152.3335
    "I34" = 1895.1666
    "J34" = 6409.5
    "K34" = 1895.1666
    "L34" = 6409.5
    "M34" = -1693.1666
    "N34" = -6813.5
    "H134" = 2094.1333
    "I134" = 1915.2307
    "K134" = 5745.6921
    "M134" = -3210.6921
    "H135" = 37779
    "J135" = 37779
    "L135" = 37779
    "N135" = -47919
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H4" = 35610650
    "I4" = 41994100
    "J4" = 24546000
    "K4" = 125982300
    "L4" = 73638000
    "M4" = -125982188
    "N4" = -73638224
    "H10" = 1521.75
    "I10" = 321.5
    "J10" = 2722
    "K10" = 964.5
    "L10" = 8166
    "M10" = -825.5
    "N10" = -8444
    "H34" = 553.0769
    "I34" = 136.25
    "J34" = 1220
    "K34" = 408.75
    "L34" = 3660
    "M34" = -324.75
    "N34" = -3828
    "H39" = 2392.5715
    "J39" = 2999.6
    "L39" = 8998.799999999999
    "N39" = -9586.799999999999
    "H55" = 1802080
    "J55" = 5000
    "L55" = 15000
    "N55" = -15354
    "H57" = 0
    "I57" = 0
    "K57" = 0
    "H58" = 4680.8
    "I58" = 405
    "J58" = 5749.75
    "K58" = 1215
    "L58" = 17249.25
    "M58" = -1087
    "N58" = -17505.25
    "H104" = 6997.857
    "J104" = 6997.857
    "L104" = 20993.571
    "N104" = -26235.571
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H106" = 37000
    "J106" = 37000
    "L106" = 37000
    "N106" = -39524
    "H123" = 59497
    "J123" = 64994.332
    "L123" = 64994.332
    "N123" = -69894.33199999999
    "H132" = 1483.1666
    "I132" = 1176.75
    "J132" = 2096
    "K132" = 3530.25
    "L132" = 6288
    "M132" = -1000.25
    "N132" = -11348
    "H133" = 137998.67
    "J133" = 137998.67
    "L133" = 137998.67
    "N133" = -148118.67
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H61" = 2906.1765
    "I61" = 2586.5
    "J61" = 4398
    "K61" = 2586.5
    "L61" = 4398
    "M61" = -2384.5
    "N61" = -4802
    "H93" = 2455.7144
    "I93" = 2207.2727
    "J93" = 3366.6667
    "K93" = 2207.2727
    "L93" = 3366.6667
    "M93" = -959.2727
    "N93" = -5862.6667
    "H113" = 2906.1765
    "I113" = 2586.5
    "J113" = 4398
    "K113" = 2586.5
    "L113" = 4398
    "M113" = -416.5
    "N113" = -8738
    "H132" = 3131
    "H133" = 42663.332
    "J133" = 42663.332
    "L133" = 42663.332
    "N133" = -47723.332
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H101" = 9949.5
    "J101" = 9949.5
    "L101" = 9949.5
    "N101" = -16439.5
    "H104" = 14599.5
    "J104" = 14599.5
    "L104" = 14599.5
    "N104" = -21587.5
    "H113" = 1291
    "I113" = 1304.9
    "J113" = 1267.8334
    "K113" = 3914.7
    "L113" = 3803.5002
    "M113" = -1744.7
    "N113" = -8143.5002
    "H132" = 2796.889
    "J132" = 3538.4614
    "L132" = 10615.3842
    "N132" = -15675.3842
    "H133" = 63935.5
    "J133" = 63935.5
    "L133" = 63935.5
    "N133" = -74055.5
    "H136" = 3134.524
    "I136" = 2074.1936
    "J136" = 6122.727
    "K136" = 6222.5808
    "L136" = 18368.181
    "M136" = -3672.5808
    "N136" = -23468.181
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
